# Update "想去人数" (interested-people count) figures in column F
# for the events whose view/interest counts changed, on both the
# "展览" sheet (source data) and the "全部类型" sheet (aggregated data).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F8").Value = 17
$wsExhibition.Range("F9").Value = 8223
$wsExhibition.Range("F11").Value = 264
$wsExhibition.Range("F12").Value = 1113
$wsExhibition.Range("F14").Value = 47
$wsExhibition.Range("F17").Value = 100
$wsExhibition.Range("F20").Value = 887

# --- Sheet "全部类型" (All types, aggregated) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 17
$wsAll.Range("F11").Value = 8223
$wsAll.Range("F13").Value = 264
$wsAll.Range("F14").Value = 1113
$wsAll.Range("F16").Value = 47
$wsAll.Range("F19").Value = 100
$wsAll.Range("F22").Value = 887

$wb.Save()
